$d = $word.ActiveDocument

# 1. WSQ funding support paragraph (full sentence replacement) -- must run
#    BEFORE the standalone title replacement below, since this sentence
#    embeds the same title substring.
$d.Content.Find.Execute(
    "We are applying for WSQ funding support for this new course PWM-Security: Security Risk Analysis (Assess and Address Security Risks) according to Security Risk Analysis Assess and Address Security Risks SEC-SRM-3002-1.1 under Security Framework.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We are applying for WSQ funding support for this new course Github Foundations Certification Training according to Software Configuration ICT-DIT-3014-1.1 under Infocomm Technology Framework.", 2)

# 2. Title: course name heading
$d.Content.Find.Execute(
    "PWM-Security: Security Risk Analysis (Assess and Address Security Risks)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Github Foundations Certification Training", 2)

# 3. Performance gaps paragraph
$d.Content.Find.Execute(
    "One of the key performance gaps in the industry is the inability to effectively identify and respond to potential security threats. Many organizations struggle with recognizing tell-tale signs of suspicious activities or behaviors, which can lead to inadequate security measures. Additionally, there is often a lack of systematic procedures for conducting thorough site assessments, resulting in vulnerabilities that could be exploited by malicious actors. This gap in risk identification and assessment can compromise the safety and security of both personnel and assets.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the rapidly evolving landscape, many organizations struggle with inefficient software release cycles and version control. Coordination issues between developers, testers, and operations teams often lead to delays and integration conflicts. Legacy systems and a lack of standardized procedures further complicate the process, increasing the risk of errors and security vulnerabilities.", 2)

# 4. "Why this course will address training needs" paragraph
$d.Content.Find.Execute(
    "This course is designed to enhance the skills necessary for identifying potential security threats through structured risk identification techniques. By focusing on systematic assessment methods, participants will learn to recognize and analyze suspicious behaviors and activities, thereby improving their ability to mitigate risks effectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This course provides hands-on experience with industry-standard tools and techniques for streamlining software development. It equips individuals with the skills to effectively manage code changes, automate deployment processes, and ensure platform-specific compatibility. By learning how to tailor software products and processes, professionals can contribute to smoother releases and improved software functionality across various platforms.", 2)

# 5. Date update
$d.Content.Find.Execute(
    "23 February 2025",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "25 February 2025", 2)
